$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.767.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -5.75%  "
$ws.Range("D3").Value = "'2.972.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.29%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'565.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.08%  "
$ws.Range("D6").Value = "'123.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.09%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").Value = "'2.965.64"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.42%  "
$ws.Range("E9").Value = "  -2.39%  "
$ws.Range("D10").Value = "'0.130"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.85%  "
$ws.Range("D11").Value = "'4.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.86%  "
$ws.Range("D12").Value = "'0.440"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.20%  "
$ws.Range("D13").Value = "'0.0000218"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.89%  "
$ws.Range("D14").Value = "'32.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.59%  "
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "'3.472.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.98%  "
$ws.Range("D17").Value = "'2.981.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.86%  "
$ws.Range("D18").Value = "'59.831.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.58%  "
$ws.Range("D19").Value = "'6.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.93%  "
$ws.Range("D20").Value = "'424.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.03%  "
$ws.Range("D21").Value = "'13.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.55%  "
$ws.Range("D22").Value = "'0.666"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.59%  "
$ws.Range("D23").Value = "'6.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.97%  "
$ws.Range("D24").Value = "'12.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.61%  "
$ws.Range("D25").Value = "'79.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.01%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").Value = "'2.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.37%  "
$ws.Range("E29").Value = "  -7.06%  "
$ws.Range("D30").Value = "'7.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.71%  "
$ws.Range("D31").Value = "'6.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.91%  "
$ws.Range("D32").Value = "'25.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.11%  "
$ws.Range("D33").Value = "'0.0964"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.43%  "
$ws.Range("D34").Value = "'5.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.87%  "
$ws.Range("E37").Value = "  -18.93%  "
$ws.Range("D38").Value = "'8.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.84%  "
$ws.Range("D39").Value = "'0.0₃0646"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -11.38%  "
$ws.Range("E40").Value = "  -9.64%  "
$ws.Range("D41").Value = "'0.106"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.97%  "
$ws.Range("D44").Value = "'2.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.88%  "
$ws.Range("D46").Value = "'120.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.01%  "
$ws.Range("D47").Value = "'0.231"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.53%  "
$ws.Range("E48").Value = "  -3.77%  "
$ws.Range("D49").Value = "'1.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.86%  "
$ws.Range("D50").Value = "'23.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.21%  "
$ws.Range("E51").Value = "  -8.28%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'50.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.32%  "
$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").Value = "'0.919"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.51%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'2.656.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.03%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'366.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.90%  "
